$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new header cell K1 (copy header styling from J1) ---
$ws.Range("K1").Value = "REGULARIZAÇÕES"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# --- Row 2: updates to existing cells ---
$ws.Range("A2").Value = "popo"
$ws.Range("C2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("G2").Value = 23
$ws.Range("H2").Value = "25/07"
$ws.Range("I2").Value = 20000
$ws.Range("J2").Value = "Pix"

# --- Row 3: new row ---
$ws.Range("A3").Value = "heloo"
$ws.Range("F3").Value = "dontknow"
$ws.Range("H3").Value = "25/05/2025"
$ws.Range("I3").Value = 500000
$ws.Range("J3").Value = "Dinheiro"

# --- Row 4: new row ---
$ws.Range("A4").Value = "Kenza Asana Rama"
$ws.Range("B4").Value = "nan"
$ws.Range("C4").Value = "nan"
$ws.Range("E4").Value = "nan"
$ws.Range("F4").Value = "nan"
$ws.Range("G4").Value = "26/07/2024"
$ws.Range("H4").Value = "26/07/2025"
# Leading apostrophe forces text storage so "35000" isn't coerced to a number,
# then reset the style so no stray quote-prefix formatting is left behind.
$ws.Range("I4").Value = "'35000"
$ws.Range("I4").Style = "Normal"
$ws.Range("J4").Value = "Pix"
$ws.Range("K4").Value = "nan"
